# Weekly refresh of the "Bruselas (repollito)" subset: the Fecha/Volumen/
# Precio columns for rows 2-27 get reshuffled across rows (a new weekly pull
# re-sorted the underlying rows). Capture a snapshot of the affected columns
# first, then write the permuted values back so reads never see partially
# updated data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that move together per source row.
$cols = @("D", "J", "K", "L", "M", "P")

# Map: destination row -> source row (values at the source row, in the
# snapshot taken below, become the new values at the destination row).
$rowMap = @{
    2  = 23
    3  = 25
    4  = 9
    5  = 5
    6  = 20
    7  = 2
    8  = 22
    9  = 7
    10 = 26
    11 = 24
    12 = 13
    13 = 21
    14 = 18
    15 = 19
    16 = 10
    17 = 3
    18 = 16
    19 = 27
    20 = 6
    21 = 11
    22 = 4
    23 = 8
    24 = 15
    25 = 12
    26 = 14
    27 = 17
}

# Snapshot current values before writing anything. Value2 (not Value) is
# used for the read so we get the raw number back instead of a wrapped
# Variant.
$snapshot = @{}
foreach ($col in $cols) {
    for ($row = 2; $row -le 27; $row++) {
        $snapshot["$col$row"] = $ws.Range("$col$row").Value2
    }
}

# Write back the permuted values.
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value2 = $snapshot["$col$srcRow"]
    }
}
